# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.075.99'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '3.317.21'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.27'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.75'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +3.14%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '3.316.98'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '3.893.55'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("D15").Value = '66.156.47'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.16'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").Value = '3.310.75'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '425.09'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -3.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.37'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.65'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -2.55%  '
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '3.461.13'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.511'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("E28").Value = '  +4.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000113'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("E32").Value = '  -2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.36'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("E36").Value = '  -3.19%  '
$ws.Range("E37").Value = '  -4.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.73'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -1.89%  '
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("D40").Value = '2.869.22'
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.35'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -3.71%  '
$ws.Range("E43").Value = '  -4.85%  '
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.83'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("E46").Value = '  -1.15%  '
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.12'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -5.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '312.59'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -2.92%  '
$ws.Range("E51").Value = '  -1.15%  '

Write-Host "Applied 73 cell updates."
